# Remove needless imports on org.dozer
#
# The "Data Mapping mappings1" Sheet1 has an import list in column C/D
# (rows 7-12) that lists classes/packages to import. Row 12 contains the
# needless "org.dozer" import which is being removed. Deleting the whole
# row shifts everything below it up by one and lets the now-unused
# "org.dozer" shared string drop out of the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Make Sheet1 the active sheet (it becomes the active tab after the edit).
$ws1.Activate() | Out-Null

# Delete the entire row that holds the "org.dozer" import (D12).
$ws1.Rows(12).Delete() | Out-Null

# Leave the selection on the cell where the deleted row used to be.
$ws1.Range("D12").Select() | Out-Null
